$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.172.48'
$ws.Range("E2").Value = '  -1.20%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.659.13'
$ws.Range("E3").Value = '  -0.87%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.004'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.80'
$ws.Range("E5").Value = '  -1.40%  '

$ws.Range("E6").Value = '  -3.07%  '

$ws.Range("E7").Value = '  +0.20%  '

$ws.Range("E8").Value = '  -1.77%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06262'
$ws.Range("E9").Value = '  -1.83%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.75'
$ws.Range("E10").Value = '  -4.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07753'
$ws.Range("E11").Value = '  -0.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.680.96'
$ws.Range("E12").Value = '  -0.53%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.444'
$ws.Range("E13").Value = '  -1.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.886.81'
$ws.Range("E14").Value = '  -0.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5431'
$ws.Range("E15").Value = '  -2.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8088'
$ws.Range("E16").Value = '  -2.96%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.87'
$ws.Range("E17").Value = '  -1.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.188.98'
$ws.Range("E18").Value = '  -1.20%  '

$ws.Range("E19").Value = '  +0.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.623'
$ws.Range("E20").Value = '  -2.95%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '192.40'
$ws.Range("E21").Value = '  +0.06%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.06'
$ws.Range("E22").Value = '  -2.57%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.008'
$ws.Range("E23").Value = '  -4.75%  '

$ws.Range("E24").Value = '  +0.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '139.99'
$ws.Range("E25").Value = '  +1.51%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1221'
$ws.Range("E26").Value = '  -4.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.224'
$ws.Range("E27").Value = '  -2.49%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.17'
$ws.Range("E28").Value = '  -0.41%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.433'
$ws.Range("E29").Value = '  +0.49%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05937'
$ws.Range("E30").Value = '  -5.45%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.271'
$ws.Range("E31").Value = '  -1.35%  '

$ws.Range("E32").Value = '  -1.16%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.256'
$ws.Range("E33").Value = '  -4.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.598'
$ws.Range("E34").Value = '  -5.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9640'

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.424'
$ws.Range("E36").Value = '  +0.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.771'
$ws.Range("E37").Value = '  -0.35%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5643'
$ws.Range("E38").Value = '  -7.99%  '

$ws.Range("E39").Value = '  -1.76%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.957'
$ws.Range("E40").Value = '  -2.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8570'
$ws.Range("E41").Value = '  -0.50%  '

$ws.Range("E42").Value = '  +0.25%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.012.28'
$ws.Range("E43").Value = '  -7.50%  '

$ws.Range("E44").Value = '  -0.51%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.801.42'
$ws.Range("E45").Value = '  -1.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈110'
$ws.Range("E46").Value = '  -2.25%  '

$ws.Range("E47").Value = '  -3.49%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.006'
$ws.Range("E48").Value = '  +0.40%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.017'
$ws.Range("E49").Value = '  -2.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05167'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.449'
$ws.Range("E51").Value = '  -5.27%  '
